$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.954.84'
$ws.Range("E2").Value = '  -2.61%  '
$ws.Range("D3").Value = '2.827.27'
$ws.Range("E3").Value = '  -2.98%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '502.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.42%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  -5.57%  '
$ws.Range("D9").Value = '2.832.26'
$ws.Range("E9").Value = '  -2.83%  '
$ws.Range("E10").Value = '  -6.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.89'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("E12").Value = '  -2.87%  '
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = '3.323.50'
$ws.Range("E14").Value = '  -3.15%  '
$ws.Range("D15").Value = '59.160.81'
$ws.Range("E15").Value = '  -2.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.04%  '
$ws.Range("D17").Value = '2.823.63'
$ws.Range("E17").Value = '  -3.63%  '
$ws.Range("E18").Value = '  -5.66%  '
$ws.Range("E19").Value = '  -6.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '347.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.60%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("E26").Value = '  -6.54%  '
$ws.Range("E27").Value = '  -7.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.31%  '
$ws.Range("E30").Value = '  -9.38%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  -4.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.900'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("D40").Value = '2.227.24'
$ws.Range("E40").Value = '  -5.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.627'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("E42").Value = '  -6.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0556'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.44%  '
$ws.Range("E45").Value = '  -10.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.21%  '
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0223'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.00%  '
$ws.Range("E49").Value = '  -4.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.99%  '
$ws.Range("E51").Value = '  -7.84%  '
